$d = $word.ActiveDocument

# The document currently has a single section with no header. Add a
# "default" header containing the questionnaire number, centered, in
# Arial 12pt, using the built-in "Header" paragraph style.

$section = $d.Sections(1)
$header = $section.Headers(1)   # wdHeaderFooterPrimary

# Insert the text (InsertAfter only materializes the single "default"
# header part instead of also minting even-page / first-page siblings).
$header.Range.InsertAfter("Questionnaire 14")

# Paragraph-level formatting: built-in Header style, centered.
$header.Range.Paragraphs(1).Style = "Header"
$header.Range.ParagraphFormat.Alignment = 1  # wdAlignParagraphCenter

# Character formatting on just the inserted text (exclude the trailing
# paragraph mark so no stray rPr ends up on the pPr).
$textRange = $header.Range.Duplicate
[void]$textRange.MoveEnd(1, -1)  # wdCharacter
$textRange.Font.Name = "Arial"
$textRange.Font.Size = 12
